# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1063
$ws1.Range("F4").Value  = 179
$ws1.Range("F5").Value  = 2907
$ws1.Range("F7").Value  = 279
$ws1.Range("F8").Value  = 27
$ws1.Range("F9").Value  = 5
$ws1.Range("F10").Value = 127
$ws1.Range("F11").Value = 100
$ws1.Range("F12").Value = 143
$ws1.Range("F13").Value = 63
$ws1.Range("F14").Value = 2738
$ws1.Range("F15").Value = 999

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1063
$ws4.Range("F5").Value  = 179
$ws4.Range("F6").Value  = 2907
$ws4.Range("F8").Value  = 279
$ws4.Range("F9").Value  = 27
$ws4.Range("F11").Value = 5
$ws4.Range("F12").Value = 127
$ws4.Range("F13").Value = 100
$ws4.Range("F14").Value = 143
$ws4.Range("F15").Value = 63
$ws4.Range("F16").Value = 2738
$ws4.Range("F17").Value = 999
